# Applies the scheduled-runner market-data refresh to the profit-tracking sheets.
# For each affected leve row, updates price/profit columns (H-N) to the refreshed
# market figures; a few rows also gain/lose a column where HQ/NQ pricing applies.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1939.963
$ws.Range("J17").Value = 1939.963
$ws.Range("L17").Value = 5819.889
$ws.Range("N17").Value = -6155.889
$ws.Range("H86").Value = 4567
$ws.Range("I86").Value = 1749.75
$ws.Range("K86").Value = 1749.75
$ws.Range("M86").Value = -626.75
$ws.Range("H89").Value = 4567
$ws.Range("I89").Value = 1749.75
$ws.Range("K89").Value = 8748.75
$ws.Range("M89").Value = -3132.75
$ws.Range("H132").Value = 3055
$ws.Range("I132").Value = 3005.3157
$ws.Range("K132").Value = 9015.947100000001
$ws.Range("M132").Value = -6485.947100000001
$ws.Range("H138").Value = 2189.1135
$ws.Range("I138").Value = 1615.3793
$ws.Range("J138").Value = 3298.3333
$ws.Range("K138").Value = 4846.1379
$ws.Range("L138").Value = 9894.999899999999
$ws.Range("M138").Value = 293.8621000000003
$ws.Range("N138").Value = -20174.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5564.3516
$ws.Range("I61").Value = 5946.778
$ws.Range("K61").Value = 5946.778
$ws.Range("M61").Value = -5734.778
$ws.Range("H74").Value = 4988.027
$ws.Range("I74").Value = 4616.7812
$ws.Range("J74").Value = 7364
$ws.Range("K74").Value = 4616.7812
$ws.Range("L74").Value = 7364
$ws.Range("M74").Value = -3742.7812
$ws.Range("N74").Value = -9112
$ws.Range("H77").Value = 4988.027
$ws.Range("I77").Value = 4616.7812
$ws.Range("J77").Value = 7364
$ws.Range("K77").Value = 23083.906
$ws.Range("L77").Value = 36820
$ws.Range("M77").Value = -18715.906
$ws.Range("N77").Value = -45556
$ws.Range("H110").Value = 725.05
$ws.Range("I110").Value = 756
$ws.Range("K110").Value = 756
$ws.Range("M110").Value = 1289
$ws.Range("H132").Value = 7638.1665
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 7638.1665
$ws.Range("K132").Value = 0
$ws.Range("N132").Value = -27974.4995
$ws.Range("H136").Value = 5564.3516
$ws.Range("I136").Value = 5946.778
$ws.Range("K136").Value = 17840.334
$ws.Range("M136").Value = -15290.334
$ws.Range("L132").Value = 22914.4995
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1304.3077
$ws.Range("I86").Value = 1087
$ws.Range("K86").Value = 1087
$ws.Range("M86").Value = 36
$ws.Range("H89").Value = 1304.3077
$ws.Range("I89").Value = 1087
$ws.Range("K89").Value = 5435
$ws.Range("M89").Value = 181
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 277.6111
$ws.Range("I7").Value = 61.5
$ws.Range("J7").Value = 547.75
$ws.Range("K7").Value = 61.5
$ws.Range("L7").Value = 547.75
$ws.Range("M7").Value = 51.5
$ws.Range("N7").Value = -773.75
$ws.Range("H11").Value = 10445
$ws.Range("J11").Value = 11519.167
$ws.Range("L11").Value = 11519.167
$ws.Range("N11").Value = -11799.167
$ws.Range("H15").Value = 1711.7778
$ws.Range("J15").Value = 7675
$ws.Range("L15").Value = 7675
$ws.Range("N15").Value = -8015
$ws.Range("H58").Value = 2720.4285
$ws.Range("J58").Value = 8007
$ws.Range("L58").Value = 8007
$ws.Range("N58").Value = -8413
$ws.Range("H86").Value = 6300
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("H89").Value = 6300
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("H132").Value = 5708.7085
$ws.Range("I132").Value = 6027.091
$ws.Range("K132").Value = 18081.273
$ws.Range("M132").Value = -15551.273
$ws.Range("H136").Value = 2720.4285
$ws.Range("J136").Value = 8007
$ws.Range("L136").Value = 24021
$ws.Range("N136").Value = -29121
$ws.Range("M86").ClearContents()
$ws.Range("M89").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1413.4286
$ws.Range("I131").Value = 603.2857
$ws.Range("J131").Value = 2223.5715
$ws.Range("K131").Value = 1809.8571
$ws.Range("L131").Value = 6670.7145
$ws.Range("M131").Value = 3230.1429
$ws.Range("N131").Value = -16750.7145

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 21329.8
$ws.Range("I132").Value = 15557
$ws.Range("K132").Value = 46671
$ws.Range("M132").Value = -44141

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4012.25
$ws.Range("I7").Value = 3854.8
$ws.Range("J7").Value = 4799.5
$ws.Range("K7").Value = 3854.8
$ws.Range("L7").Value = 4799.5
$ws.Range("M7").Value = -3742.8
$ws.Range("N7").Value = -5023.5
$ws.Range("H22").Value = 144912.86
$ws.Range("I22").Value = 250600
$ws.Range("J22").Value = 3996.6667
$ws.Range("K22").Value = 250600
$ws.Range("L22").Value = 3996.6667
$ws.Range("M22").Value = -250305
$ws.Range("N22").Value = -4586.6667
$ws.Range("H27").Value = 144912.86
$ws.Range("I27").Value = 250600
$ws.Range("J27").Value = 3996.6667
$ws.Range("K27").Value = 250600
$ws.Range("L27").Value = 3996.6667
$ws.Range("M27").Value = -250493
$ws.Range("N27").Value = -4210.6667
$ws.Range("H43").Value = 756250
$ws.Range("J43").Value = 756250
$ws.Range("L43").Value = 756250
$ws.Range("N43").Value = -756636
$ws.Range("H61").Value = 5469.9
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 5469.9
$ws.Range("K61").Value = 0
$ws.Range("N61").Value = -5873.9
$ws.Range("H107").Value = 19398.8
$ws.Range("I107").Value = 19398.8
$ws.Range("K107").Value = 19398.8
$ws.Range("M107").Value = -17478.8
$ws.Range("H113").Value = 5469.9
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 5469.9
$ws.Range("K113").Value = 0
$ws.Range("N113").Value = -9809.9
$ws.Range("H122").Value = 3779.1082
$ws.Range("I122").Value = 3530.484
$ws.Range("K122").Value = 10591.452
$ws.Range("M122").Value = -8141.451999999999
$ws.Range("H126").Value = 4012.25
$ws.Range("I126").Value = 3854.8
$ws.Range("J126").Value = 4799.5
$ws.Range("K126").Value = 11564.4
$ws.Range("L126").Value = 14398.5
$ws.Range("M126").Value = -9094.400000000001
$ws.Range("N126").Value = -19338.5
$ws.Range("H132").Value = 3736
$ws.Range("I132").Value = 3400.1333
$ws.Range("J132").Value = 4995.5
$ws.Range("K132").Value = 10200.3999
$ws.Range("L132").Value = 14986.5
$ws.Range("M132").Value = -7670.3999
$ws.Range("N132").Value = -20046.5
$ws.Range("L61").Value = 5469.9
$ws.Range("L113").Value = 5469.9
$ws.Range("M61").ClearContents()
$ws.Range("M113").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 96490.2
$ws.Range("J62").Value = 7237.875
$ws.Range("L62").Value = 7237.875
$ws.Range("N62").Value = -8485.875
$ws.Range("H65").Value = 96490.2
$ws.Range("J65").Value = 7237.875
$ws.Range("L65").Value = 36189.375
$ws.Range("N65").Value = -42429.375
$ws.Range("H81").Value = 1010.5455
$ws.Range("I81").Value = 888
$ws.Range("J81").Value = 1225
$ws.Range("K81").Value = 1776
$ws.Range("L81").Value = 2450
$ws.Range("M81").Value = -715
$ws.Range("N81").Value = -4572
$ws.Range("H84").Value = 1010.5455
$ws.Range("I84").Value = 888
$ws.Range("J84").Value = 1225
$ws.Range("K84").Value = 8880
$ws.Range("L84").Value = 12250
$ws.Range("M84").Value = -3576
$ws.Range("N84").Value = -22858
$ws.Range("H96").Value = 1151.2
$ws.Range("J96").Value = 849.8
$ws.Range("L96").Value = 849.8
$ws.Range("N96").Value = -3595.8
$ws.Range("H107").Value = 2439.9583
$ws.Range("I107").Value = 1640.5385
$ws.Range("J107").Value = 3384.7273
$ws.Range("K107").Value = 4921.6155
$ws.Range("L107").Value = 10154.1819
$ws.Range("M107").Value = -3001.6155
$ws.Range("N107").Value = -13994.1819
$ws.Range("H122").Value = 2066.7693
$ws.Range("I122").Value = 1640.6957
$ws.Range("K122").Value = 4922.0871
$ws.Range("M122").Value = -2472.0871
$ws.Range("H132").Value = 7678.2915
$ws.Range("I132").Value = 6654.1113
$ws.Range("K132").Value = 19962.3339
$ws.Range("M132").Value = -17432.3339
